$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 217 - this shifts the existing rows 217..353
# down to 218..354, carrying their values/formatting with them.
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new weekly price record.
$ws.Cells.Item(217, 1).Value = 5
$ws.Cells.Item(217, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(217, 3).Value = "Maule"
$ws.Cells.Item(217, 4).Value = 44767
$ws.Cells.Item(217, 5).Value = 7
$ws.Cells.Item(217, 6).Value = 100114014
$ws.Cells.Item(217, 7).Value = "Betarraga"
$ws.Cells.Item(217, 8).Value = "Sin especificar"
$ws.Cells.Item(217, 9).Value = "Primera"
$ws.Cells.Item(217, 10).Value = 5000
$ws.Cells.Item(217, 11).Value = 750
$ws.Cells.Item(217, 12).Value = 750
$ws.Cells.Item(217, 13).Value = 750
$ws.Cells.Item(217, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(217, 15).Value = "Región del Maule"
$ws.Cells.Item(217, 16).Value = 150
$ws.Cells.Item(217, 17).Value = 5
$ws.Cells.Item(217, 18).Value = "Hortaliza"
